# Corrección modificar cita vista cliente
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 23/24 status swap:
#    F23: PROGRESS (yellow) -> DONE (green)
#    F24: FLUTTER (blue)    -> PROGRESS (yellow)
$ws.Range("F23").Value = "DONE"
$ws.Range("F23").Interior.Color = 5287936   # FF00B050 green, same as F2 DONE style

$ws.Range("F24").Value = "PROGRESS"
$ws.Range("F24").Interior.Color = 65535     # FFFFFF00 yellow, same as old F23 PROGRESS style

# 2) Move the TODO note up from H36 to H33, updating its text, keeping the red fill style
$note = "PASAR A PROD APPLICATIONPROPERTIES"
$ws.Range("H33").Value = $note
$ws.Range("H33").Interior.Color = 192       # FFC00000 dark red

# 3) Clear the old H36 note cell (and F36/G36) to a plain, no-fill style
$ws.Range("F36:H36").ClearContents()
$ws.Range("F36:H36").Interior.ColorIndex = -4142   # xlColorIndexNone

# 4) Update the selection shown when the sheet is opened
$ws.Range("F24").Select()
